$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 790; this shifts old rows 790-831 down to 791-832
# and old row 789 (2026/02/06, 11, 201) stays put.
$ws.Rows.Item(790).Insert()

$ws.Range("A790").NumberFormat = "@"
$ws.Range("A790").Value = "2026/02/06"
$ws.Range("A790").Style = "Normal"
$ws.Range("B790").Value = "金"
$ws.Range("C790").Value = 14
$ws.Range("D790").Value = 201

$ws.Range("A791").NumberFormat = "@"
$ws.Range("A791").Value = "2026/12/29"
$ws.Range("A791").Style = "Normal"
$ws.Range("B791").Value = "火"
$ws.Range("C791").Value = 13
$ws.Range("D791").Value = 88

$ws.Range("A792").NumberFormat = "@"
$ws.Range("A792").Value = "2026/12/29"
$ws.Range("A792").Style = "Normal"
$ws.Range("B792").Value = "火"
$ws.Range("C792").Value = 16
$ws.Range("D792").Value = 99

$ws.Range("A793").NumberFormat = "@"
$ws.Range("A793").Value = "2026/12/29"
$ws.Range("A793").Style = "Normal"
$ws.Range("B793").Value = "火"
$ws.Range("C793").Value = 19
$ws.Range("D793").Value = 81

$ws.Range("A794").NumberFormat = "@"
$ws.Range("A794").Value = "2026/12/29"
$ws.Range("A794").Style = "Normal"
$ws.Range("B794").Value = "火"
$ws.Range("C794").Value = 23
$ws.Range("D794").Value = 85

$ws.Range("A795").NumberFormat = "@"
$ws.Range("A795").Value = "2026/12/30"
$ws.Range("A795").Style = "Normal"
$ws.Range("B795").Value = "水"
$ws.Range("C795").Value = 2
$ws.Range("D795").Value = 89

$ws.Range("A796").NumberFormat = "@"
$ws.Range("A796").Value = "2026/12/30"
$ws.Range("A796").Style = "Normal"
$ws.Range("B796").Value = "水"
$ws.Range("C796").Value = 5
$ws.Range("D796").Value = 95

$ws.Range("A797").NumberFormat = "@"
$ws.Range("A797").Value = "2026/12/30"
$ws.Range("A797").Style = "Normal"
$ws.Range("B797").Value = "水"
$ws.Range("C797").Value = 8
$ws.Range("D797").Value = 91

$ws.Range("A798").NumberFormat = "@"
$ws.Range("A798").Value = "2026/12/30"
$ws.Range("A798").Style = "Normal"
$ws.Range("B798").Value = "水"
$ws.Range("C798").Value = 13
$ws.Range("D798").Value = 92

$ws.Range("A799").NumberFormat = "@"
$ws.Range("A799").Value = "2026/12/30"
$ws.Range("A799").Style = "Normal"
$ws.Range("B799").Value = "水"
$ws.Range("C799").Value = 16
$ws.Range("D799").Value = 99

$ws.Range("A800").NumberFormat = "@"
$ws.Range("A800").Value = "2026/12/30"
$ws.Range("A800").Style = "Normal"
$ws.Range("B800").Value = "水"
$ws.Range("C800").Value = 22
$ws.Range("D800").Value = 108

$ws.Range("A801").NumberFormat = "@"
$ws.Range("A801").Value = "2026/12/31"
$ws.Range("A801").Style = "Normal"
$ws.Range("B801").Value = "木"
$ws.Range("C801").Value = 2
$ws.Range("D801").Value = 114

$ws.Range("A802").NumberFormat = "@"
$ws.Range("A802").Value = "2026/12/31"
$ws.Range("A802").Style = "Normal"
$ws.Range("B802").Value = "木"
$ws.Range("C802").Value = 6
$ws.Range("D802").Value = 120

$ws.Range("A803").NumberFormat = "@"
$ws.Range("A803").Value = "2026/12/31"
$ws.Range("A803").Style = "Normal"
$ws.Range("B803").Value = "木"
$ws.Range("C803").Value = 9
$ws.Range("D803").Value = 120

$ws.Range("A804").NumberFormat = "@"
$ws.Range("A804").Value = "2026/12/31"
$ws.Range("A804").Style = "Normal"
$ws.Range("B804").Value = "木"
$ws.Range("C804").Value = 12
$ws.Range("D804").Value = 201

$ws.Range("A805").NumberFormat = "@"
$ws.Range("A805").Value = "2026/12/31"
$ws.Range("A805").Style = "Normal"
$ws.Range("B805").Value = "木"
$ws.Range("C805").Value = 14
$ws.Range("D805").Value = 130

$ws.Range("A806").NumberFormat = "@"
$ws.Range("A806").Value = "2026/12/31"
$ws.Range("A806").Style = "Normal"
$ws.Range("B806").Value = "木"
$ws.Range("C806").Value = 22
$ws.Range("D806").Value = 120

$ws.Range("A807").NumberFormat = "@"
$ws.Range("A807").Value = "2027/01/01"
$ws.Range("A807").Style = "Normal"
$ws.Range("B807").Value = "金"
$ws.Range("C807").Value = 2
$ws.Range("D807").Value = 129

$ws.Range("A808").NumberFormat = "@"
$ws.Range("A808").Value = "2027/01/01"
$ws.Range("A808").Style = "Normal"
$ws.Range("B808").Value = "金"
$ws.Range("C808").Value = 5
$ws.Range("D808").Value = 119

$ws.Range("A809").NumberFormat = "@"
$ws.Range("A809").Value = "2027/01/01"
$ws.Range("A809").Style = "Normal"
$ws.Range("B809").Value = "金"
$ws.Range("C809").Value = 13
$ws.Range("D809").Value = 133

$ws.Range("A810").NumberFormat = "@"
$ws.Range("A810").Value = "2027/01/01"
$ws.Range("A810").Style = "Normal"
$ws.Range("B810").Value = "金"
$ws.Range("C810").Value = 16
$ws.Range("D810").Value = 109

$ws.Range("A811").NumberFormat = "@"
$ws.Range("A811").Value = "2027/01/01"
$ws.Range("A811").Style = "Normal"
$ws.Range("B811").Value = "金"
$ws.Range("C811").Value = 19
$ws.Range("D811").Value = 120

$ws.Range("A812").NumberFormat = "@"
$ws.Range("A812").Value = "2027/01/02"
$ws.Range("A812").Style = "Normal"
$ws.Range("B812").Value = "土"
$ws.Range("C812").Value = 1
$ws.Range("D812").Value = 105

$ws.Range("A813").NumberFormat = "@"
$ws.Range("A813").Value = "2027/01/02"
$ws.Range("A813").Style = "Normal"
$ws.Range("B813").Value = "土"
$ws.Range("C813").Value = 5
$ws.Range("D813").Value = 109

$ws.Range("A814").NumberFormat = "@"
$ws.Range("A814").Value = "2027/01/02"
$ws.Range("A814").Style = "Normal"
$ws.Range("B814").Value = "土"
$ws.Range("C814").Value = 8
$ws.Range("D814").Value = 110

$ws.Range("A815").NumberFormat = "@"
$ws.Range("A815").Value = "2027/01/02"
$ws.Range("A815").Style = "Normal"
$ws.Range("B815").Value = "土"
$ws.Range("C815").Value = 13
$ws.Range("D815").Value = 132

$ws.Range("A816").NumberFormat = "@"
$ws.Range("A816").Value = "2027/01/02"
$ws.Range("A816").Style = "Normal"
$ws.Range("B816").Value = "土"
$ws.Range("C816").Value = 16
$ws.Range("D816").Value = 145

$ws.Range("A817").NumberFormat = "@"
$ws.Range("A817").Value = "2027/01/02"
$ws.Range("A817").Style = "Normal"
$ws.Range("B817").Value = "土"
$ws.Range("C817").Value = 19
$ws.Range("D817").Value = 157

$ws.Range("A818").NumberFormat = "@"
$ws.Range("A818").Value = "2027/01/02"
$ws.Range("A818").Style = "Normal"
$ws.Range("B818").Value = "土"
$ws.Range("C818").Value = 22
$ws.Range("D818").Value = 165

$ws.Range("A819").NumberFormat = "@"
$ws.Range("A819").Value = "2027/01/03"
$ws.Range("A819").Style = "Normal"
$ws.Range("B819").Value = "日"
$ws.Range("C819").Value = 1
$ws.Range("D819").Value = 174

$ws.Range("A820").NumberFormat = "@"
$ws.Range("A820").Value = "2027/01/03"
$ws.Range("A820").Style = "Normal"
$ws.Range("B820").Value = "日"
$ws.Range("C820").Value = 4
$ws.Range("D820").Value = 192

$ws.Range("A821").NumberFormat = "@"
$ws.Range("A821").Value = "2027/01/03"
$ws.Range("A821").Style = "Normal"
$ws.Range("B821").Value = "日"
$ws.Range("C821").Value = 7
$ws.Range("D821").Value = 189

$ws.Range("A822").NumberFormat = "@"
$ws.Range("A822").Value = "2027/01/03"
$ws.Range("A822").Style = "Normal"
$ws.Range("B822").Value = "日"
$ws.Range("C822").Value = 13
$ws.Range("D822").Value = 201

$ws.Range("A823").NumberFormat = "@"
$ws.Range("A823").Value = "2027/01/03"
$ws.Range("A823").Style = "Normal"
$ws.Range("B823").Value = "日"
$ws.Range("C823").Value = 16
$ws.Range("D823").Value = 201

$ws.Range("A824").NumberFormat = "@"
$ws.Range("A824").Value = "2027/01/03"
$ws.Range("A824").Style = "Normal"
$ws.Range("B824").Value = "日"
$ws.Range("C824").Value = 19
$ws.Range("D824").Value = 201

$ws.Range("A825").NumberFormat = "@"
$ws.Range("A825").Value = "2027/01/03"
$ws.Range("A825").Style = "Normal"
$ws.Range("B825").Value = "日"
$ws.Range("C825").Value = 22
$ws.Range("D825").Value = 194

$ws.Range("A826").NumberFormat = "@"
$ws.Range("A826").Value = "2027/01/04"
$ws.Range("A826").Style = "Normal"
$ws.Range("B826").Value = "月"
$ws.Range("C826").Value = 2
$ws.Range("D826").Value = 164

$ws.Range("A827").NumberFormat = "@"
$ws.Range("A827").Value = "2027/01/04"
$ws.Range("A827").Style = "Normal"
$ws.Range("B827").Value = "月"
$ws.Range("C827").Value = 5
$ws.Range("D827").Value = 166

$ws.Range("A828").NumberFormat = "@"
$ws.Range("A828").Value = "2027/01/04"
$ws.Range("A828").Style = "Normal"
$ws.Range("B828").Value = "月"
$ws.Range("C828").Value = 7
$ws.Range("D828").Value = 168

$ws.Range("A829").NumberFormat = "@"
$ws.Range("A829").Value = "2027/01/04"
$ws.Range("A829").Style = "Normal"
$ws.Range("B829").Value = "月"
$ws.Range("C829").Value = 13
$ws.Range("D829").Value = 173

$ws.Range("A830").NumberFormat = "@"
$ws.Range("A830").Value = "2027/01/04"
$ws.Range("A830").Style = "Normal"
$ws.Range("B830").Value = "月"
$ws.Range("C830").Value = 22
$ws.Range("D830").Value = 127

$ws.Range("A831").NumberFormat = "@"
$ws.Range("A831").Value = "2027/01/05"
$ws.Range("A831").Style = "Normal"
$ws.Range("B831").Value = "火"
$ws.Range("C831").Value = 1
$ws.Range("D831").Value = 118

$ws.Range("A832").NumberFormat = "@"
$ws.Range("A832").Value = "2027/01/05"
$ws.Range("A832").Style = "Normal"
$ws.Range("B832").Value = "火"
$ws.Range("C832").Value = 7
$ws.Range("D832").Value = 127
